$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.110.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.403.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -9.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -2.92%  "

$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.25"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("E12").Value = "  -2.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.96"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.948.43"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.83%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.503.13"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.19%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.119"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.38%  "

$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.162.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.985"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "408.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.73"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.72"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "588.35"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("E33").Value = "  -4.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.66"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.95%  "

$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "35.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.99%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.196.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.35%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.370"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.03%  "

$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("E41").Value = "  -11.52%  "

$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("E44").Value = "  -5.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.37%  "

$ws.Range("E46").Value = "  -3.97%  "

$ws.Range("E47").Value = "  -5.95%  "

$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.66%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.97%  "
